# Auto-generated Excel COM-interop script
# Updates market price / profit figures (columns H-N) for several leve rows
# across multiple crafting-class worksheets (scheduled market-data refresh).

$wb = $excel.ActiveWorkbook

# --- Worksheet "ALC" ---
$ws = $wb.Worksheets.Item("ALC")
# Row 10
$ws.Range("H10").Value = 30000
$ws.Range("I10").Value = 30000
$ws.Range("K10").Value = 30000
$ws.Range("M10").Value = -29707
# Row 43
$ws.Range("H43").Value = 999.8570999999999
$ws.Range("J43").Value = 1066.3334
$ws.Range("L43").Value = 1066.3334
$ws.Range("N43").Value = -1204.3334
# Row 98
$ws.Range("H98").Value = 733.069
$ws.Range("J98").Value = 555
$ws.Range("L98").Value = 555
$ws.Range("N98").Value = -3551
# Row 112
$ws.Range("H112").Value = 2092.739
$ws.Range("J112").Value = 1917.579
$ws.Range("L112").Value = 5752.737
$ws.Range("N112").Value = -7968.737
# Row 122
$ws.Range("H122").Value = 733.069
$ws.Range("J122").Value = 555
$ws.Range("L122").Value = 1665
$ws.Range("N122").Value = -6565
# Row 127
$ws.Range("H127").Value = 1297.8
$ws.Range("I127").Value = 1228.2858
$ws.Range("K127").Value = 3684.8574
$ws.Range("M127").Value = 1275.1426
# Row 132
$ws.Range("H132").Value = 3796.6758
$ws.Range("I132").Value = 3430.7188
$ws.Range("K132").Value = 10292.1564
$ws.Range("M132").Value = -7762.1564
# Row 135
$ws.Range("H135").Value = 955.3913
$ws.Range("I135").Value = 652.8421
$ws.Range("J135").Value = 2392.5
$ws.Range("K135").Value = 5875.5789
$ws.Range("L135").Value = 21532.5
$ws.Range("M135").Value = -3340.5789
$ws.Range("N135").Value = -26602.5
# Row 138
$ws.Range("H138").Value = 1925.4894
$ws.Range("I138").Value = 1656.1111
$ws.Range("J138").Value = 2289.15
$ws.Range("K138").Value = 4968.3333
$ws.Range("L138").Value = 6867.450000000001
$ws.Range("M138").Value = 171.6666999999998
$ws.Range("N138").Value = -17147.45

# --- Worksheet "ARM" ---
$ws = $wb.Worksheets.Item("ARM")
# Row 33
$ws.Range("H33").Value = 8188.727
$ws.Range("I33").Value = 8188.727
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 8188.727
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -7859.727
$ws.Range("N33").Value = ""

# --- Worksheet "CRP" ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5264.061
$ws.Range("I31").Value = 10411.077
$ws.Range("J31").Value = 3405.4167
$ws.Range("K31").Value = 10411.077
$ws.Range("L31").Value = 3405.4167
$ws.Range("M31").Value = -10116.077
$ws.Range("N31").Value = -3995.4167
# Row 34
$ws.Range("H34").Value = 5264.061
$ws.Range("I34").Value = 10411.077
$ws.Range("J34").Value = 3405.4167
$ws.Range("K34").Value = 10411.077
$ws.Range("L34").Value = 3405.4167
$ws.Range("M34").Value = -10209.077
$ws.Range("N34").Value = -3809.4167
# Row 58
$ws.Range("H58").Value = 1799.4242
$ws.Range("I58").Value = 1566.8695
$ws.Range("J58").Value = 2334.3
$ws.Range("K58").Value = 1566.8695
$ws.Range("L58").Value = 2334.3
$ws.Range("M58").Value = -1363.8695
$ws.Range("N58").Value = -2740.3
# Row 132
$ws.Range("H132").Value = 2908.2
$ws.Range("I132").Value = 1519
$ws.Range("K132").Value = 4557
$ws.Range("M132").Value = -2027
# Row 136
$ws.Range("H136").Value = 1799.4242
$ws.Range("I136").Value = 1566.8695
$ws.Range("J136").Value = 2334.3
$ws.Range("K136").Value = 4700.6085
$ws.Range("L136").Value = 7002.900000000001
$ws.Range("M136").Value = -2150.6085
$ws.Range("N136").Value = -12102.9

# --- Worksheet "CUL" ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 11773299
$ws.Range("I4").Value = 20179860
$ws.Range("K4").Value = 60539580
$ws.Range("M4").Value = -60539468
# Row 82
$ws.Range("H82").Value = 14603.182
$ws.Range("J82").Value = 15012.857
$ws.Range("L82").Value = 45038.571
$ws.Range("N82").Value = -45850.571
# Row 85
$ws.Range("H85").Value = 14603.182
$ws.Range("J85").Value = 15012.857
$ws.Range("L85").Value = 45038.571
$ws.Range("N85").Value = -47846.571
# Row 129
$ws.Range("H129").Value = 2132.5386
$ws.Range("I129").Value = 1264
$ws.Range("J129").Value = 2877
$ws.Range("K129").Value = 3792
$ws.Range("L129").Value = 8631
$ws.Range("M129").Value = 1208
$ws.Range("N129").Value = -18631

# --- Worksheet "GSM" ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3893.3845
$ws.Range("J80").Value = 3950
$ws.Range("L80").Value = 3950
$ws.Range("N80").Value = -5946
# Row 83
$ws.Range("H83").Value = 3893.3845
$ws.Range("J83").Value = 3950
$ws.Range("L83").Value = 19750
$ws.Range("N83").Value = -29734
# Row 132
$ws.Range("H132").Value = 24110.773
$ws.Range("I132").Value = 35497.848
$ws.Range("J132").Value = 7662.778
$ws.Range("K132").Value = 106493.544
$ws.Range("L132").Value = 22988.334
$ws.Range("M132").Value = -103963.544
$ws.Range("N132").Value = -28048.334

# --- Worksheet "LTW" ---
$ws = $wb.Worksheets.Item("LTW")
# Row 43
$ws.Range("H43").Value = 1710000
$ws.Range("J43").Value = 2515000
$ws.Range("L43").Value = 2515000
$ws.Range("N43").Value = -2515386
# Row 63
$ws.Range("H63").Value = 46687.25
$ws.Range("J63").Value = 47583
$ws.Range("L63").Value = 47583
$ws.Range("N63").Value = -49081
# Row 66
$ws.Range("H66").Value = 46687.25
$ws.Range("J66").Value = 47583
$ws.Range("L66").Value = 142749
$ws.Range("N66").Value = -150237
# Row 82
$ws.Range("H82").Value = 3008.4546
$ws.Range("I82").Value = 3811.75
$ws.Range("J82").Value = 866.3333
$ws.Range("K82").Value = 3811.75
$ws.Range("L82").Value = 866.3333
$ws.Range("M82").Value = -3450.75
$ws.Range("N82").Value = -1588.3333
# Row 85
$ws.Range("H85").Value = 3008.4546
$ws.Range("I85").Value = 3811.75
$ws.Range("J85").Value = 866.3333
$ws.Range("K85").Value = 3811.75
$ws.Range("L85").Value = 866.3333
$ws.Range("M85").Value = -2563.75
$ws.Range("N85").Value = -3362.3333
# Row 132
$ws.Range("H132").Value = 6440.9565
$ws.Range("I132").Value = 3591.077
$ws.Range("J132").Value = 10145.8
$ws.Range("K132").Value = 10773.231
$ws.Range("L132").Value = 30437.4
$ws.Range("M132").Value = -8243.231
$ws.Range("N132").Value = -35497.39999999999
# Row 136
$ws.Range("H136").Value = 1946.5682
$ws.Range("I136").Value = 1617.6
$ws.Range("J136").Value = 3225.889
$ws.Range("K136").Value = 4852.799999999999
$ws.Range("L136").Value = 9677.667000000001
$ws.Range("M136").Value = -2302.799999999999
$ws.Range("N136").Value = -14777.667

# --- Worksheet "WVR" ---
$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 4997
$ws.Range("I14").Value = 3746.25
$ws.Range("K14").Value = 3746.25
$ws.Range("M14").Value = -3578.25
# Row 122
$ws.Range("H122").Value = 1783.3103
$ws.Range("J122").Value = 2184.1428
$ws.Range("L122").Value = 6552.428400000001
$ws.Range("N122").Value = -11452.4284
# Row 136
$ws.Range("H136").Value = 3102.15
$ws.Range("I136").Value = 931.9286
$ws.Range("K136").Value = 2795.7858
$ws.Range("M136").Value = -245.7857999999997

